$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits inside the
#    "Modificar Servlet" paragraph of the "reseña" section (Sergio's
#    cell). We find the paragraph that owns the bookmark and rewrite it
#    without the bookmarkStart/bookmarkEnd pair, leaving everything else
#    (the text + the spell-check proofErr markers) untouched.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm.Start

$paras = $d.Paragraphs
$oldPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    $candRange = $cand.Range
    if ($candRange.Start -le $bmStart -and $bmStart -le $candRange.End) {
        $oldPara = $candRange
        break
    }
}

$noBookmarkXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:ind w:left="0"/>
              <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Modificar </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Servlet</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$oldPara.InsertXML($noBookmarkXml) | Out-Null

# ---------------------------------------------------------------------
# 2) Fill in the (until now empty) third-column cell belonging to
#    "Miguel David Yac Castro - 2022020" with the list of items he
#    worked on, finishing with a "Modificar Servlet" line that now owns
#    the "_GoBack" bookmark instead.
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)
$targetRow = 0
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    if ($table.Cell($r, 1).Range.Text -like "*2022020*") {
        $targetRow = $r
        break
    }
}

$targetCell = $table.Cell($targetRow, 3)

$cellXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:ind w:left="0"/>
              <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t>Modelo Proveedor</w:t>
            </w:r>
            <w:r>
              <w:t>es</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:ind w:left="0"/>
              <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t>Modelo DAO Proveedor</w:t>
            </w:r>
            <w:r>
              <w:t>es</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:ind w:left="0"/>
              <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t>Vista Proveedor</w:t>
            </w:r>
            <w:r>
              <w:t>es</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:ind w:left="0"/>
              <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
              <w:rPr>
                <w:u w:val="single"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:t>Controlador Proveedor</w:t>
            </w:r>
            <w:r>
              <w:t>es</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:ind w:left="0"/>
              <w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Modificar </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Servlet</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$targetCell.Range.InsertXML($cellXml) | Out-Null

Write-Output "Edit applied successfully"
